$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.126.96"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.13%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.627.67"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "600.64"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "146.00"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.05%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.627.41"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  -0.68%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.59"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  +0.00%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.361"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.49%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.09"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.08%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.098.43"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.25%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "63.003.30"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("E17").Value = "  -2.21%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.626.13"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("E20").Value = "  +2.17%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "339.52"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.86"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.87%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -3.57%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "66.42"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("E26").Value = "  -3.47%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.62"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.48%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "547.57"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("E29").Value = "  -7.85%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.162"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.97%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.83"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("E35").Value = "  -2.39%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +10.63%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "165.75"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -0.08%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.93"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("E41").Value = "  +5.73%  "
$ws.Range("E42").Value = "  +0.03%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "167.52"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.27%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.72"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.41%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "22.29"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  +0.30%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0958"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.52"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "
